$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: p_houseGasBurnerEfficiency_r
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "p_houseGasBurnerEfficiency_r"
$ws.Range("C7").Value = 0.95
$ws.Range("E7").Value = "Residential gas burner efficiency"

# Row 9: p_houseHeatPumpEfficiency_r (shared strings entered before row 8's, to
# match the original authoring/shared-string insertion order)
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "p_houseHeatPumpEfficiency_r"
$ws.Range("C9").Value = 3
$ws.Range("E9").Value = "Residential heatPump efficiency (COP = 3)"

# Row 8: p_houseHeatPumpElectricCapacity_kW
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "p_houseHeatPumpElectricCapacity_kW"
$ws.Range("C8").Value = 11
$ws.Range("D8").Value = "kW"
$ws.Range("E8").Value = "Average residential heatpump electrical power"

# Update selection to match diff (C17 selected, even though data only to row 9)
$ws.Range("C17").Select()
